$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the numeric bounds used by the batch test row (row 2) from
# 4000/20000 down to 2000/10000, and derived counts from 16000 down to 8000.
$ws.Range("I2").Value = "10000"
$ws.Range("L2").Value = "10000"
$ws.Range("M2").Value = "select count(*) from `$schema26 where id>2000 and id<=10000"
$ws.Range("N2").Value = "8000"
$ws.Range("O2").Value = "update `$schema26 set name='BJ' where id>2000 and id<=10000"
$ws.Range("P2").Value = "8000"
$ws.Range("R2").Value = "8000"
$ws.Range("T2").Value = "10000"

# Move the active selection from C7 to C8, as recorded in the sheet view.
$ws.Range("C8").Select()
